$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 2-11, columns B:G
$data = @(
    @(2,  -0.2243426211202774, 2.382314769375817,  19.09061844583581, 4.369281227597487, 4.46620300197608,  22),
    @(3,  -0.7120868115500602, 1.933618923903567,  10.88254014902561, 3.298869525917267, 3.300643293046703, 21),
    @(4,  -0.3051420100610381, 1.578682252322016,  10.2251943824045,  3.197685785439917, 3.265784782946779, 20),
    @(5,  -0.2804117665781,    1.683207935177004,  11.04924039551323, 3.324039770446983, 3.402952828034286, 19),
    @(6,  -0.3173356620072107, 1.884213172660337,  10.94775440033029, 3.308739095234057, 3.388969517186058, 18),
    @(7,  -0.2356007239738259, 1.940837229566887,  12.01945225291491, 3.466908169091721, 3.565345857874665, 17),
    @(8,  -0.1538742921933148, 1.930586848953809,  10.98523007623265, 3.314397392623983, 3.41940389010477,  16),
    @(9,  -0.1230455590706185, 1.907399787585887,  12.46868919605704, 3.531103113200893, 3.652819212576486, 15),
    @(10, -0.09178185018074302,2.13619364861378,   13.88641155430669, 3.726447578365579, 3.865944277603689, 14),
    @(11, -0.08394991111291034,2.334233754498058,  14.90545946209773, 3.860758923074286, 4.017455189220993, 13)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
